$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove obsolete city rows (London1, Paris2, Madrid$$$) and the trailing
# blank row, leaving just City / Amsterdan / Rome
$ws.Rows.Item(2).Delete()   # London1 (row shifts Amsterdan up to row 2)
$ws.Rows.Item(3).Delete()   # Paris2 (row shifts Rome up to row 3)
$ws.Rows.Item(4).Delete()   # Madrid$$$
$ws.Rows.Item(4).Delete()   # trailing blank row

# Autofit column A to match content width and move the active selection
$ws.Columns.Item(1).AutoFit()
$ws.Range("F5").Select()
